$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.000" or "22.024.21"
# are not reinterpreted as numbers/dates by Excel's smart-parsing.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '22.024.21'
$ws.Range("E2").Value = '  -1.84%  '

# Row 3
$ws.Range("D3").Value = '1.554.40'
$ws.Range("E3").Value = '  -1.19%  '

# Row 4
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("E5").Value = '  -0.03%  '

# Row 6
$ws.Range("D6").Value = '286.61'
$ws.Range("E6").Value = '  -0.43%  '

# Row 7
$ws.Range("D7").Value = '0.3769'
$ws.Range("E7").Value = '  +1.64%  '

# Row 8
$ws.Range("D8").Value = '0.3237'
$ws.Range("E8").Value = '  -2.37%  '

# Row 9
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '41.15'
$ws.Range("E9").Value = '  -12.65%  '

# Row 10
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.122'
$ws.Range("E10").Value = '  -2.93%  '

# Row 11
$ws.Range("D11").Value = '0.07299'
$ws.Range("E11").Value = '  -2.71%  '

# Row 12
$ws.Range("E12").Value = '  -0.06%  '

# Row 13
$ws.Range("D13").Value = '19.29'
$ws.Range("E13").Value = '  -6.94%  '

# Row 14
$ws.Range("D14").Value = '5.712'
$ws.Range("E14").Value = '  -3.71%  '

# Row 15
$ws.Range("E15").Value = '  -1.66%  '

# Row 16
$ws.Range("D16").Value = '1.552.08'
$ws.Range("E16").Value = '  -0.67%  '

# Row 17
$ws.Range("D17").Value = '0.00001078'
$ws.Range("E17").Value = '  -3.23%  '

# Row 18
$ws.Range("D18").Value = '0.06648'
$ws.Range("E18").Value = '  -1.17%  '

# Row 19
$ws.Range("D19").Value = '85.00'
$ws.Range("E19").Value = '  -3.76%  '

# Row 20
$ws.Range("D20").Value = '6.422'
$ws.Range("E20").Value = '  +0.25%  '

# Row 21
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").Value = '  -0.05%  '

# Row 22
$ws.Range("E22").Value = '  -3.25%  '

# Row 23
$ws.Range("E23").Value = '  -3.84%  '

# Row 24
$ws.Range("D24").Value = '22.041.70'
$ws.Range("E24").Value = '  -1.73%  '

# Row 25
$ws.Range("D25").Value = '2.263'
$ws.Range("E25").Value = '  -4.64%  '

# Row 26
$ws.Range("D26").Value = '2.508'
$ws.Range("E26").Value = '  -4.45%  '

# Row 27
$ws.Range("D27").Value = '149.86'
$ws.Range("E27").Value = '  -0.66%  '

# Row 28
$ws.Range("D28").Value = '18.85'
$ws.Range("E28").Value = '  -3.60%  '

# Row 29
$ws.Range("D29").Value = '4.851'
$ws.Range("E29").Value = '  -1.93%  '

# Row 30
$ws.Range("D30").Value = '1.729.16'
$ws.Range("E30").Value = '  -0.61%  '

# Row 31
$ws.Range("D31").Value = '120.05'
$ws.Range("E31").Value = '  -3.93%  '

# Row 32
$ws.Range("D32").Value = '1.120'
$ws.Range("E32").Value = '  +2.68%  '

# Row 33
$ws.Range("D33").Value = '5.912'
$ws.Range("E33").Value = '  -2.72%  '

# Row 34
$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.649'
$ws.Range("E34").Value = '  -17.04%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.08164'
$ws.Range("E35").Value = '  -1.87%  '

# Row 36
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '9.270'
$ws.Range("E36").Value = '  -6.02%  '

# Row 37
$ws.Range("D37").Value = '5.230'
$ws.Range("E37").Value = '  -1.79%  '

# Row 38
$ws.Range("D38").Value = '0.02281'

# Row 39
$ws.Range("E39").Value = '  -3.54%  '

# Row 40
$ws.Range("E40").Value = '  -4.76%  '

# Row 41
$ws.Range("D41").Value = '1.212'
$ws.Range("E41").Value = '  -7.04%  '

# Row 42
$ws.Range("E42").Value = '  -4.52%  '

# Row 44
$ws.Range("D44").Value = '0.5923'
$ws.Range("E44").Value = '  -4.95%  '

# Row 45
$ws.Range("D45").Value = '13.43'
$ws.Range("E45").Value = '  -3.80%  '

# Row 46
$ws.Range("D46").Value = '3.723'
$ws.Range("E46").Value = '  -1.33%  '

# Row 47
$ws.Range("D47").Value = '0.5725'
$ws.Range("E47").Value = '  -5.19%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '120.00'
$ws.Range("E48").Value = '  -3.83%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.929'
$ws.Range("E49").Value = '  -5.30%  '

# Row 50
$ws.Range("E50").Value = '  -4.31%  '

# Row 51
$ws.Range("D51").Value = '0.06897'
$ws.Range("E51").Value = '  -4.13%  '
